$d = $word.ActiveDocument

# Locate the list paragraph that holds the "Senha" / ";" runs (the one that
# also carries the _GoBack bookmark straddling them).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Senha*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Step 1: remove the "Senha" text run, scoped to this paragraph only.
    # The bookmark start/end that sits between the "Senha" run and the ";"
    # run is left untouched by the search/replace.
    $pRange = $d.Paragraphs($targetIndex).Range
    $pRange.Find.Execute("Senha", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 2) | Out-Null

    # Step 2: merge this paragraph (now just ";" plus the bookmark) with the
    # following one ("Nome;") by deleting the paragraph mark between them.
    $pRange = $d.Paragraphs($targetIndex).Range
    $paraMark = $d.Range($pRange.End - 1, $pRange.End)
    $paraMark.Delete()

    # Step 3: drop the leftover ";" that used to terminate the "Senha" item,
    # leaving the merged paragraph reading "Nome;" with the bookmark intact.
    $merged = $d.Paragraphs($targetIndex).Range
    $leadingSemicolon = $d.Range($merged.Start, $merged.Start + 1)
    $leadingSemicolon.Delete()
}
